# Fruta / hortaliza, semanal
# Insert two new weekly records (date 2023-05-31) above the existing
# row 84 block, shifting the subsequent rows (84-96) down to (86-98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 84, pushing everything from row 84 down by two rows.
$ws.Range("A84:A85").EntireRow.Insert()

# Populate the two newly inserted rows with the new weekly data.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are identical to the surrounding records.

$ws.Cells.Item(84, 1).Value = 8
$ws.Cells.Item(84, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(84, 3).Value = "Coquimbo"
$ws.Cells.Item(84, 4).Value = 45077
$ws.Cells.Item(84, 5).Value = 4
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100104
$ws.Cells.Item(84, 8).Value = "Frutos de pepita"
$ws.Cells.Item(84, 9).Value = 100104003
$ws.Cells.Item(84, 10).Value = "Membrillo"
$ws.Cells.Item(84, 11).Value = "Champion"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 14
$ws.Cells.Item(84, 14).Value = 240000
$ws.Cells.Item(84, 15).Value = 250000
$ws.Cells.Item(84, 16).Value = 245000
$ws.Cells.Item(84, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(84, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(84, 19).Value = 544
$ws.Cells.Item(84, 20).Value = 450

$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 45077
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100104
$ws.Cells.Item(85, 8).Value = "Frutos de pepita"
$ws.Cells.Item(85, 9).Value = 100104003
$ws.Cells.Item(85, 10).Value = "Membrillo"
$ws.Cells.Item(85, 11).Value = "Champion"
$ws.Cells.Item(85, 12).Value = "Segunda"
$ws.Cells.Item(85, 13).Value = 12
$ws.Cells.Item(85, 14).Value = 210000
$ws.Cells.Item(85, 15).Value = 220000
$ws.Cells.Item(85, 16).Value = 215000
$ws.Cells.Item(85, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(85, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(85, 19).Value = 478
$ws.Cells.Item(85, 20).Value = 450

# Ensure the date cells keep the same date number format used by the rest
# of column D.
$ws.Range("D84:D85").NumberFormat = $ws.Range("D86").NumberFormat
